# Actualización asociación usuario entrevista
# Rename the "Reclutador" column to "Usuario" and set its values to "login".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Postulantes")

# Header (row 1, column C)
$ws.Range("C1").Value = "Usuario"

# Data rows (rows 2-4, column C) -> all become "login"
$ws.Range("C2").Value = "login"
$ws.Range("C3").Value = "login"
$ws.Range("C4").Value = "login"

# Restore the active selection to A4, matching the saved workbook state.
$ws.Range("A4").Select()
